$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 160.9723563333334
$ws.Range("H2").Value = 482.917069
$ws.Range("I2").Value = 0.3931645655589854
$ws.Range("J2").Value = 0.3931645655589854
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.05968133333333333
$ws.Range("N2").Value = 0.179044
$ws.Range("O2").Value = 0.02602747651633847
$ws.Range("P2").Value = 0.02602747651633848
$ws.Range("Q2").Value = 9.607044855781778
$ws.Range("R2").Value = 86.46340370203599
$ws.Range("S2").Value = 0.01023308149714291
$ws.Range("T2").Value = 0.01023308149714291

$ws.Range("G3").Value = 160.9723563333334
$ws.Range("H3").Value = 482.917069
$ws.Range("I3").Value = 0.3931645655589854
$ws.Range("J3").Value = 0.3931645655589854
$ws.Range("O3").Value = 0.144012433133819
$ws.Range("P3").Value = 0.144012433133819
$ws.Range("Q3").Value = 53.15666711055812
$ws.Range("R3").Value = 478.410003995023
$ws.Range("S3").Value = 0.05662058570815038
$ws.Range("T3").Value = 0.05662058570815038

$ws.Range("G4").Value = 160.9723563333334
$ws.Range("H4").Value = 482.917069
$ws.Range("I4").Value = 0.3931645655589854
$ws.Range("J4").Value = 0.3931645655589854
$ws.Range("O4").Value = 0.8299600903498424
$ws.Range("P4").Value = 0.8299600903498425
$ws.Range("Q4").Value = 306.3479400891737
$ws.Range("R4").Value = 2757.131460802563
$ws.Range("S4").Value = 0.326310898353692
$ws.Range("T4").Value = 0.3263108983536921

$ws.Range("I5").Value = 0.2197635343237224
$ws.Range("J5").Value = 0.2197635343237224
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.05968133333333333
$ws.Range("N5").Value = 0.179044
$ws.Range("O5").Value = 0.02602747651633847
$ws.Range("P5").Value = 0.02602747651633848
$ws.Range("Q5").Value = 5.369960359757777
$ws.Range("R5").Value = 48.32964323781999
$ws.Range("S5").Value = 0.005719890228758228
$ws.Range("T5").Value = 0.005719890228758229

$ws.Range("I6").Value = 0.2197635343237224
$ws.Range("J6").Value = 0.2197635343237224
$ws.Range("O6").Value = 0.144012433133819
$ws.Range("P6").Value = 0.144012433133819
$ws.Range("S6").Value = 0.03164868129204681
$ws.Range("T6").Value = 0.03164868129204681

$ws.Range("I7").Value = 0.2197635343237224
$ws.Range("J7").Value = 0.2197635343237224
$ws.Range("O7").Value = 0.8299600903498424
$ws.Range("P7").Value = 0.8299600903498425
$ws.Range("S7").Value = 0.1823949628029173
$ws.Range("T7").Value = 0.1823949628029174

$ws.Range("I8").Value = 0.3870719001172923
$ws.Range("J8").Value = 0.3870719001172923
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.05968133333333333
$ws.Range("N8").Value = 0.179044
$ws.Range("O8").Value = 0.02602747651633847
$ws.Range("P8").Value = 0.02602747651633848
$ws.Range("Q8").Value = 9.458169511162668
$ws.Range("R8").Value = 85.12352560046399
$ws.Range("S8").Value = 0.01007450479043734
$ws.Range("T8").Value = 0.01007450479043734

$ws.Range("I9").Value = 0.3870719001172923
$ws.Range("J9").Value = 0.3870719001172923
$ws.Range("O9").Value = 0.144012433133819
$ws.Range("P9").Value = 0.144012433133819
$ws.Range("Q9").Value = 52.33292606909468
$ws.Range("R9").Value = 470.996334621852
$ws.Range("S9").Value = 0.05574316613362182
$ws.Range("T9").Value = 0.05574316613362183

$ws.Range("I10").Value = 0.3870719001172923
$ws.Range("J10").Value = 0.3870719001172923
$ws.Range("O10").Value = 0.8299600903498424
$ws.Range("P10").Value = 0.8299600903498425
$ws.Range("S10").Value = 0.321254229193233
$ws.Range("T10").Value = 0.3212542291932332

$wb.Save()
